$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093260407447815
$ws.Range("B1").Value = 1.885980606079102
$ws.Range("D1").Value = 1.240871429443359
$ws.Range("E1").Value = 1.167752027511597
